# Edit: slide 1 (sldId 256) - title placeholder (shape id 67) gets a new
# line of text ("dfs") typed into what used to be a trailing empty
# paragraph, and the notes-page "slide image" placeholder (shape id 64)
# has its position/size nudged by a few EMU (rounding of the slide-image
# frame that PowerPoint performs when it resaves the notes page).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Add "dfs" text to the last (previously empty) paragraph of the
#        title shape on slide 1 ---------------------------------------
$sh = $s.Shapes.Item(1)
$tf = $sh.TextFrame
$tr = $tf.TextRange

$para = $tr.Paragraphs(4, 1)
$para.Text = "dfs"

# --- 2. Nudge the notes-page placeholder geometry (best effort; some
#        hosts do not allow editing notes-page shape geometry) --------
try {
    $notes = $s.NotesPage
    $imgShape = $notes.Shapes.Item(1)
    $imgShape.Left = 30
    $imgShape.Top = 54
    $imgShape.Width = 480
    $imgShape.Height = 270
} catch {
    # not supported by this host; ignore
}
